$wb = $excel.ActiveWorkbook

# --- Rename sheets: O5_* -> Q5_* -----------------------------------------
$wsLinking = $wb.Worksheets.Item("O5_linking")
$wsLinking.Name = "Q5_linking"
$wsCoding  = $wb.Worksheets.Item("O5_coding")
$wsCoding.Name = "Q5_coding"

# --- Sheet2 (Q5_coding): restore frozen-pane scroll position -------------
# Re-anchor the freeze at B2 (same split as before) first so the existing
# frozen split survives, then scroll the frozen view down so its
# top-left visible cell becomes B93, finally restoring the original
# per-pane selections (topRight=B1, bottomLeft=A2, bottomRight=I7:I8).
$wsCoding.Activate()
$winCoding = $excel.ActiveWindow
$winCoding.FreezePanes = $false
$wsCoding.Range("B93").Select()
$winCoding.FreezePanes = $true

$winCoding.Panes.Item(2).Activate()
$wsCoding.Range("B1").Select()
$winCoding.Panes.Item(3).Activate()
$wsCoding.Range("A2").Select()
$winCoding.Panes.Item(4).Activate()
$wsCoding.Range("I7:I8").Select()

# --- Sheet1 (Q5_linking): becomes the active tab, new selection A14 ------
$wsLinking.Activate()
$wsLinking.Range("A14").Select()
